# Weekly fruit/vegetable price update: insert a new weekly record at the
# top of the "Cebollín baby" data block (row 58), pushing the existing
# rows 58-75 down to 59-76.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 58-75 down to 59-76, inserting a blank row 58.
$ws.Rows("58:58").Insert()

# Populate the newly inserted row 58 with this week's record. All the
# "boilerplate" columns (market/region/category/etc.) are identical to
# every other row in this block.
$ws.Range("A58").Value = 1
$ws.Range("B58").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C58").Value = "Arica y Parinacota"
$ws.Range("D58").Value = 44627
$ws.Range("E58").Value = 15
$ws.Range("F58").Value = 100112038
$ws.Range("G58").Value = "Cebollín baby"
$ws.Range("H58").Value = "Sin especificar"
$ws.Range("I58").Value = "Primera"
$ws.Range("J58").Value = 300
$ws.Range("K58").Value = 2000
$ws.Range("L58").Value = 2500
$ws.Range("M58").Value = 2250
$ws.Range("N58").Value = "$/paquete 1,5 a 2 kilos"
$ws.Range("O58").Value = "Región de Arica y Parinacota"
$ws.Range("P58").Value = 1125
$ws.Range("Q58").Value = 2
$ws.Range("R58").Value = "Hortaliza"

Write-Output "Inserted new weekly row at 58; rows shifted to 76."
